$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 10 new rows before row 11, shifting the existing rows (11-22) down to (21-32)
$ws.Rows("11:20").Insert()

# Fill the newly inserted rows in column E with the new parameter labels
$ws.Range("E11").Value = "R1"
$ws.Range("E12").Value = "R2"
$ws.Range("E13").Value = "v_aw"
$ws.Range("E14").Value = "v_alpha_aw"
$ws.Range("E15").Value = "v_ar"
$ws.Range("E16").Value = "alpha_ar"
$ws.Range("E17").Value = "v_ak"
$ws.Range("E18").Value = "alpha_ak"
$ws.Range("E19").Value = "v_ah"
$ws.Range("E20").Value = "alpha_ah"

# The row that used to hold "J" (now row 25) did not have a "Completed?" mark; add "yes"
$ws.Range("G25").Value = "yes"

# Leave the selection on the last newly-entered cell
$ws.Range("E20").Select() | Out-Null
